$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the formatting of the other
# header cells (e.g. G1: bold, bordered, centered/top-aligned).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data for the new column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
